# fix: sexting intensity escalation + PPV0 phase detection bug across all 23 models
# Update JessicaFPJourney sheet (sexting escalation script column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JessicaFPJourney")

$ws.Range("B4").Value = "cum with me amor... right now, don't look away 🥵"
$ws.Range("B6").Value = "I'm cumming for you right now... watch me, every second of it"
$ws.Range("B7").Value = "I'm on the edge amor... I've been holding back and I can't anymore, I need to let go 🥵"
$ws.Range("B8").Value = "dios mio 😏"
$ws.Range("B9").Value = "you're about to see something I don't show anyone amor... this is all because of you 🥵"
$ws.Range("B11").Value = "I'm about to give you something you'll never forget amor... watch this"
$ws.Range("B12").Value = "my fingers are going so deep and fast and I can feel every wave building inside me 🥵"
$ws.Range("B13").Value = "I'm fucking myself right now and all I can picture is your face between my legs... dios mio"
$ws.Range("B14").Value = "FUCK 😏"
$ws.Range("B15").Value = "mira lo que me haces amor... this is what you do to me 🥵"
$ws.Range("B17").Value = "tell me what you'd do to me right now amor... I want to hear every dirty detail"
$ws.Range("B18").Value = "I need to feel something inside me right now... I keep imagining it's you and it's driving me crazy 🥵"
$ws.Range("B19").Value = "I'm so wet right now I can feel it on my thighs... this is what you do to me amor"
$ws.Range("B20").Value = "dios... that was intense but I'm just getting started 😏"
$ws.Range("B21").Value = "look what you started amor... you have no idea what's coming next 🥵"
$ws.Range("B23").Value = "I know exactly what I want right now amor and it starts with you watching what I do next 🥵"
$ws.Range("B24").Value = "I'm lying here with my hand between my thighs and it's because this conversation is getting to me"
$ws.Range("B25").Value = "I could tell you'd like that... and honestly? your reaction is turning me on more than I expected 😏"
